# The header row previously read "Keywords" / "Suggestions" — update it to
# the new "Question" / "Answer" headers. All other rows keep their existing
# text; the shared-string table gets rebuilt/compacted automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Question"
$ws.Range("B1").Value = "Answer"

# Reset the view: scroll back to the top and make B1 the active/selected cell
# (previously the sheet was scrolled to A58 with B79 selected).
$ws.Range("B1").Select()
